# "primera actualizacion poblada con numeros" -- populate the per-skill
# score columns (AK..AX) that were previously empty/placeholder on most
# player report sheets, and correct the provisional numbers on "Extremo".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Central
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Central")
$ws.Range("AK3").Value = 85
$ws.Range("AL3").Value = 90
$ws.Range("AM3").Value = 75
$ws.Range("AN3").Value = 60
$ws.Range("AO3").Value = 80
$ws.Range("AP3").Value = 80
$ws.Range("AQ3").Value = 80
$ws.Range("AR3").Value = 95
$ws.Range("AS3").Value = 85
$ws.Range("AT3").Value = 80
$ws.Range("AU3").Value = 90

# ---------------------------------------------------------------------
# Lat Izq
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Lat Izq")
$ws.Range("AK2").Value = 85
$ws.Range("AL2").Value = 90
$ws.Range("AM2").Value = 75
$ws.Range("AN2").Value = 80
$ws.Range("AO2").Value = 80
$ws.Range("AP2").Value = 90
$ws.Range("AQ2").Value = 85
$ws.Range("AR2").Value = 95
$ws.Range("AS2").Value = 85
$ws.Range("AT2").Value = 80
$ws.Range("AU2").Value = 90
$ws.Range("AV2").Value = 80
$ws.Range("AW2").Value = 70
$ws.Range("AX2").Value = 80

$ws.Range("AK3").Value = 80
$ws.Range("AL3").Value = 80
$ws.Range("AM3").Value = 90
$ws.Range("AN3").Value = 85
$ws.Range("AO3").Value = 75
$ws.Range("AP3").Value = 90
$ws.Range("AQ3").Value = 70
$ws.Range("AR3").Value = 80
$ws.Range("AS3").Value = 70
$ws.Range("AT3").Value = 80
$ws.Range("AU3").Value = 85
$ws.Range("AV3").Value = 80
$ws.Range("AW3").Value = 75
$ws.Range("AX3").Value = 80

# ---------------------------------------------------------------------
# Mixto
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Mixto")
$ws.Range("AK2").Value = 85
$ws.Range("AL2").Value = 90
$ws.Range("AM2").Value = 75
$ws.Range("AN2").Value = 60
$ws.Range("AO2").Value = 80
$ws.Range("AP2").Value = 80
$ws.Range("AQ2").Value = 80
$ws.Range("AR2").Value = 95
$ws.Range("AS2").Value = 85
$ws.Range("AT2").Value = 80
$ws.Range("AU2").Value = 90
$ws.Range("AV2").Value = 70
$ws.Range("AW2").Value = 85
$ws.Range("AX2").Value = 75

# ---------------------------------------------------------------------
# Ofensivo
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Ofensivo")
$ws.Range("AK2").Value = 85
$ws.Range("AL2").Value = 90
$ws.Range("AM2").Value = 75
$ws.Range("AN2").Value = 80
$ws.Range("AO2").Value = 80
$ws.Range("AP2").Value = 90
$ws.Range("AQ2").Value = 85
$ws.Range("AR2").Value = 95
$ws.Range("AS2").Value = 85
$ws.Range("AT2").Value = 80
$ws.Range("AU2").Value = 90

$ws.Range("AK3").Value = 80
$ws.Range("AL3").Value = 80
$ws.Range("AM3").Value = 90
$ws.Range("AN3").Value = 85
$ws.Range("AO3").Value = 75
$ws.Range("AP3").Value = 90
$ws.Range("AQ3").Value = 70
$ws.Range("AR3").Value = 80
$ws.Range("AS3").Value = 70
$ws.Range("AT3").Value = 80
$ws.Range("AU3").Value = 85

# ---------------------------------------------------------------------
# Extremo -- row 2 had provisional numbers, correct them; rows 3 & 4 get
# their scores filled in for the first time.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Extremo")
$ws.Range("AL2").Value = 65
$ws.Range("AM2").Value = 80
$ws.Range("AN2").Value = 30
$ws.Range("AO2").Value = 80
$ws.Range("AP2").Value = 70
$ws.Range("AS2").Value = 70
$ws.Range("AU2").Value = 85

$ws.Range("AK3").Value = 85
$ws.Range("AL3").Value = 90
$ws.Range("AM3").Value = 75
$ws.Range("AN3").Value = 60
$ws.Range("AO3").Value = 80
$ws.Range("AP3").Value = 80
$ws.Range("AQ3").Value = 85
$ws.Range("AR3").Value = 95
$ws.Range("AS3").Value = 85
$ws.Range("AT3").Value = 80
$ws.Range("AU3").Value = 90

$ws.Range("AK4").Value = 100
$ws.Range("AL4").Value = 75
$ws.Range("AM4").Value = 90
$ws.Range("AN4").Value = 70
$ws.Range("AO4").Value = 75
$ws.Range("AP4").Value = 75
$ws.Range("AQ4").Value = 45
$ws.Range("AR4").Value = 95
$ws.Range("AS4").Value = 80
$ws.Range("AT4").Value = 80
$ws.Range("AU4").Value = 85

$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 26
$ws.Range("AK2:AU4").Select()

# ---------------------------------------------------------------------
# Centrodelantero
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Centrodelantero")
$ws.Range("AK2").Value = 90
$ws.Range("AL2").Value = 65
$ws.Range("AM2").Value = 80
$ws.Range("AN2").Value = 60
$ws.Range("AO2").Value = 80
$ws.Range("AP2").Value = 70
$ws.Range("AQ2").Value = 70
$ws.Range("AR2").Value = 89
$ws.Range("AS2").Value = 70
$ws.Range("AT2").Value = 75
$ws.Range("AU2").Value = 85
$ws.Range("AV2").Value = 75

$ws.Range("AK3").Value = 85
$ws.Range("AL3").Value = 90
$ws.Range("AM3").Value = 75
$ws.Range("AN3").Value = 60
$ws.Range("AO3").Value = 80
$ws.Range("AP3").Value = 80
$ws.Range("AQ3").Value = 80
$ws.Range("AR3").Value = 95
$ws.Range("AS3").Value = 85
$ws.Range("AT3").Value = 80
$ws.Range("AU3").Value = 90
$ws.Range("AV3").Value = 90

$ws.Range("AK4").Value = 90
$ws.Range("AL4").Value = 75
$ws.Range("AM4").Value = 90
$ws.Range("AN4").Value = 70
$ws.Range("AO4").Value = 90
$ws.Range("AP4").Value = 75
$ws.Range("AQ4").Value = 90
$ws.Range("AR4").Value = 95
$ws.Range("AS4").Value = 80
$ws.Range("AT4").Value = 80
$ws.Range("AU4").Value = 85
$ws.Range("AV4").Value = 75

# ---------------------------------------------------------------------
# Update each sheet's on-screen selection / scroll position to match
# where the author was working, then leave "Centrodelantero" active
# (it was, and remains, the selected tab).
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Central")
$ws.Activate()
$ws.Range("A6:AK8").Select()

$ws = $wb.Worksheets.Item("Lat Izq")
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 23
$ws.Range("AK2:AY3").Select()

$ws = $wb.Worksheets.Item("Mixto")
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 29
$ws.Range("AI30").Select()

$ws = $wb.Worksheets.Item("Ofensivo")
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 25
$ws.Range("AK2:AU3").Select()

$ws = $wb.Worksheets.Item("Centrodelantero")
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 34
$ws.Range("AV3").Select()
